{"js": "// The commit adds an explanatory sentence to the end of the\n// \"\uc5c5\uce90\uc2a4\ud305\uacfc \ub2e4\uc6b4\uce90\uc2a4\ud305\" (upcasting/downcasting) bullet, turning it into a\n// \"term : definition\" line like its sibling bullets (e.g. \"\uc0dd\uc131\uc790 : ...\",\n// \"\ud074\ub798\uc2a4 \uc0c1\uc18d : ...\", \"\uba54\uc11c\ub4dc \uc624\ubc84\ub77c\uc774\ub529 : ...\").\n//\n// Final text should read:\n//   \uc5c5\uce90\uc2a4\ud305\uacfc \ub2e4\uc6b4\uce90\uc2a4\ud305 : \uc8fc\uc18c\uac12\uc73c\ub85c \ud560\ub2f9\ud558\ub294 \uacbd\uc6b0\ub9cc \uc801\uc6a9\ub418\uace0 \uc2a4\ud0dd\uc73c\ub85c \ub4e4\uc5b4\uac00\uba74 \uadf8\ub0e5 \ud615\ubcc0\ud658\uc774 \uc77c\uc5b4\ub098\uac8c \ub41c\ub2e4.\n\nconst searchText = \"\uc5c5\uce90\uc2a4\ud305\uacfc \ub2e4\uc6b4\uce90\uc2a4\ud305\";\nconst appendText =\n  \" : \uc8fc\uc18c\uac12\uc73c\ub85c \ud560\ub2f9\ud558\ub294 \uacbd\uc6b0\ub9cc \uc801\uc6a9\ub418\uace0 \uc2a4\ud0dd\uc73c\ub85c \ub4e4\uc5b4\uac00\uba74 \uadf8\ub0e5 \ud615\ubcc0\ud658\uc774 \uc77c\uc5b4\ub098\uac8c \ub41c\ub2e4.\";\n\nconst results = context.document.body.search(searchText, {\n  matchCase: true,\n  matchWholeWord: false,\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find paragraph containing \"' + searchText + '\"');\n}\n\n// Use the first match; get its paragraph and append the new sentence at the\n// very end of it, mirroring the run/formatting pattern Word uses elsewhere\n// in this document for mixed Korean/English \"term : definition\" bullets\n// (Korean segments get rFonts hint=\"eastAsia\"; \"\uc8fc\uc18c\uac12\uc73c\ub85c\" is wrapped the\n// way Word's proofer brackets a single spell-checked word).\nconst target = results.items[0];\nconst paragraph = target.paragraphs.getFirst();\n\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>' +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  \"</Relationships>\" +\n  \"</pkg:xmlData></pkg:part>\" +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p>' +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">: </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t>\uc8fc\uc18c\uac12\uc73c\ub85c</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/></w:rPr><w:t xml:space=\"preserve\"> \ud560\ub2f9\ud558\ub294 \uacbd\uc6b0\ub9cc \uc801\uc6a9\ub418\uace0 \uc2a4\ud0dd\uc73c\ub85c \ub4e4\uc5b4\uac00\uba74 \uadf8\ub0e5 \ud615\ubcc0\ud658\uc774 \uc77c\uc5b4\ub098\uac8c \ub41c\ub2e4.</w:t></w:r>' +\n  \"</w:p></w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part>\" +\n  \"</pkg:package>\";\n\ntry {\n  // Preferred: insert real OOXML so run boundaries/formatting match Word's\n  // own output exactly.\n  paragraph.insertOoxml(ooxml, Word.InsertLocation.end);\n  await context.sync();\n} catch (e) {\n  // Fallback: plain-text append still yields the correct visible content\n  // even if OOXML insertion isn't supported by the host.\n  paragraph.insertText(appendText, Word.InsertLocation.end);\n  await context.sync();\n}\n", "ps1": "# The commit adds an explanatory sentence to the end of the\n# \"\uc5c5\uce90\uc2a4\ud305\uacfc \ub2e4\uc6b4\uce90\uc2a4\ud305\" (upcasting/downcasting) bullet, turning it into a\n# \"term : definition\" line like its sibling bullets (e.g. \"\uc0dd\uc131\uc790 : ...\",\n# \"\ud074\ub798\uc2a4 \uc0c1\uc18d : ...\", \"\uba54\uc11c\ub4dc \uc624\ubc84\ub77c\uc774\ub529 : ...\").\n#\n# Final text should read:\n#   \uc5c5\uce90\uc2a4\ud305\uacfc \ub2e4\uc6b4\uce90\uc2a4\ud305 : \uc8fc\uc18c\uac12\uc73c\ub85c \ud560\ub2f9\ud558\ub294 \uacbd\uc6b0\ub9cc \uc801\uc6a9\ub418\uace0 \uc2a4\ud0dd\uc73c\ub85c \ub4e4\uc5b4\uac00\uba74 \uadf8\ub0e5 \ud615\ubcc0\ud658\uc774 \uc77c\uc5b4\ub098\uac8c \ub41c\ub2e4.\n\n$d = $word.ActiveDocument\n\n$searchText = \"\uc5c5\uce90\uc2a4\ud305\uacfc \ub2e4\uc6b4\uce90\uc2a4\ud305\"\n$appendText = \" : \uc8fc\uc18c\uac12\uc73c\ub85c \ud560\ub2f9\ud558\ub294 \uacbd\uc6b0\ub9cc \uc801\uc6a9\ub418\uace0 \uc2a4\ud0dd\uc73c\ub85c \ub4e4\uc5b4\uac00\uba74 \uadf8\ub0e5 \ud615\ubcc0\ud658\uc774 \uc77c\uc5b4\ub098\uac8c \ub41c\ub2e4.\"\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = $searchText\n$find.Forward = $true\n$find.Wrap = 0\n\n$found = $find.Execute()\n\nif ($found) {\n    # Grow the found range to cover the whole paragraph, then drop the\n    # trailing paragraph mark so the insertion lands at the end of the\n    # visible text (not at the start of the following paragraph).\n    [void]$range.Expand(4)            # wdParagraph\n    [void]$range.MoveEnd(1, -1)       # wdCharacter, shrink by 1 (paragraph mark)\n    $range.InsertAfter($appendText)\n}\n"}
